# Add team record columns (Wins/Losses/Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new labels in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header styling (bold, centered, thin box border)
$hdr = $ws.Range("AD1:AF1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.LineStyle = 1         # xlContinuous

# Data rows 2-54: same team record repeated for every player
$ws.Range("AD2:AD54").Value = 95
$ws.Range("AE2:AE54").Value = 67
$ws.Range("AF2:AF54").Value = 0
